# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells should carry the same bold/centered/bordered style as
# the existing header row (style index 1), so copy formatting from an
# existing header cell before setting the text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every data row (2-52) gets the team's season record.
$ws.Range("AD2:AD52").Value = 76
$ws.Range("AE2:AE52").Value = 86
$ws.Range("AF2:AF52").Value = 0
